# Data retrieved - Tue Jun  8 18:52:08 UTC 2021
#
# Apply the tiny re-save precision fix to the last existing row's
# timestamp (A41) and append the newly retrieved row (row 42) of job
# numbers data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A41's timestamp is re-serialized with slightly different float precision
# on this re-save (same instant, just a different trailing digit).
$ws.Cells.Item(41, 1).Value2 = 44354.8317763206

# Append the new data row (row 42) retrieved on 2021-06-08.
$ws.Cells.Item(42, 1).Value2 = 44355.78620761175
$ws.Cells.Item(42, 2).Value2 = 75731
$ws.Cells.Item(42, 3).Value2 = 63840
$ws.Cells.Item(42, 4).Value2 = 3404
$ws.Cells.Item(42, 5).Value2 = 2101
$ws.Cells.Item(42, 6).Value2 = 1484
$ws.Cells.Item(42, 7).Value2 = 19917
$ws.Cells.Item(42, 8).Value2 = 1433
$ws.Cells.Item(42, 9).Value2 = 882
$ws.Cells.Item(42, 10).Value2 = 203
